# Update "想去人数" (want-to-go count, column F) for specific rows on two
# sheets: "展览" and "全部类型". These are refreshed scrape counters — the
# row/column layout and every other cell stay untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 56
$ws1.Range("F3").Value = 3221
$ws1.Range("F5").Value = 2272
$ws1.Range("F6").Value = 331
$ws1.Range("F7").Value = 321
$ws1.Range("F8").Value = 1199
$ws1.Range("F9").Value = 1038
$ws1.Range("F10").Value = 257
$ws1.Range("F11").Value = 474
$ws1.Range("F14").Value = 79
$ws1.Range("F16").Value = 8034
$ws1.Range("F17").Value = 350
$ws1.Range("F18").Value = 2471
$ws1.Range("F23").Value = 546
$ws1.Range("F27").Value = 1853
$ws1.Range("F28").Value = 290
$ws1.Range("F30").Value = 1684
$ws1.Range("F32").Value = 1909
$ws1.Range("F34").Value = 3
$ws1.Range("F35").Value = 54
$ws1.Range("F37").Value = 285
$ws1.Range("F38").Value = 44
$ws1.Range("F39").Value = 191
$ws1.Range("F40").Value = 358
$ws1.Range("F42").Value = 228

# --- Sheet: 全部类型 -------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 56
$ws4.Range("F5").Value = 3221
$ws4.Range("F7").Value = 2272
$ws4.Range("F8").Value = 331
$ws4.Range("F9").Value = 321
$ws4.Range("F10").Value = 1199
$ws4.Range("F12").Value = 1038
$ws4.Range("F13").Value = 257
$ws4.Range("F14").Value = 474
$ws4.Range("F16").Value = 79
$ws4.Range("F18").Value = 8035
$ws4.Range("F19").Value = 350
$ws4.Range("F20").Value = 2471
$ws4.Range("F26").Value = 546
$ws4.Range("F30").Value = 1853
$ws4.Range("F31").Value = 291
$ws4.Range("F33").Value = 1684
$ws4.Range("F35").Value = 1909
$ws4.Range("F37").Value = 3
$ws4.Range("F38").Value = 54
$ws4.Range("F40").Value = 285
$ws4.Range("F41").Value = 44
$ws4.Range("F42").Value = 191
$ws4.Range("F43").Value = 358
$ws4.Range("F49").Value = 228
